# Append the "Friday, Jan 13" departures to the RZE_Departures log.
# Columns: A=NUMBER  B=DATE  C=TIME  D=FLIGHT  E=TO  F=SHORT
#          G=AIRLINE H=MODEL I=AIRCFAT ID J=STATUS K=(blank) L=DIFFERENCE M=(blank)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row=53; Number=52; Date="Friday, Jan 13"; Time="4:45 AM";  Flight="AEG482"; To="Frankfurt"; Short="(FRA)"; Airline="Airest ";             Model="SF34"; AircraftId="(ES-LSI)"; Status="4:41 AM";  Difference="0 hours, -4 minutes" }
    @{ Row=54; Number=53; Date="Friday, Jan 13"; Time="5:40 AM";  Flight="LO3804"; To="Warsaw";    Short="(WAW)"; Airline="LOT ";                Model="E190"; AircraftId="(SP-LMC)"; Status="5:53 AM";  Difference="0 hours, 13 minutes" }
    @{ Row=55; Number=54; Date="Friday, Jan 13"; Time="8:00 AM";  Flight="X7592";  To="Liege";      Short="(LGG)"; Airline="Challenge Airlines "; Model="B744"; AircraftId="(OO-ACE)"; Status="9:39 AM";  Difference="1 hours, 39 minutes" }
    @{ Row=56; Number=55; Date="Friday, Jan 13"; Time="12:00 PM"; Flight="SK7182"; To="Berlin";     Short="(BER)"; Airline="SAS ";                Model="B737"; AircraftId="(LN-RPJ)"; Status="12:42 PM"; Difference="0 hours, 42 minutes" }
    @{ Row=57; Number=56; Date="Friday, Jan 13"; Time="12:15 PM"; Flight="LO3810"; To="Warsaw";     Short="(WAW)"; Airline="LOT ";                Model="E195"; AircraftId="(SP-LNK)"; Status="12:15 PM"; Difference="0 hours, 0 minutes" }
    @{ Row=58; Number=57; Date="Friday, Jan 13"; Time="1:35 PM";  Flight="FR8225"; To="Bristol";    Short="(BRS)"; Airline="Ryanair ";            Model="B738"; AircraftId="(EI-DCM)"; Status="1:50 PM";  Difference="0 hours, 15 minutes" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Number
    $ws.Range("B$row").Value = $r.Date
    $ws.Range("C$row").Value = $r.Time
    $ws.Range("D$row").Value = $r.Flight
    $ws.Range("E$row").Value = $r.To
    $ws.Range("F$row").Value = $r.Short
    $ws.Range("G$row").Value = $r.Airline
    $ws.Range("H$row").Value = $r.Model
    $ws.Range("I$row").Value = $r.AircraftId
    $ws.Range("J$row").Value = $r.Status
    $ws.Range("L$row").Value = $r.Difference
}
